# Natmi following Dr Hou advice
# Rewrite the LR-pairs result table (sheet1) for Efna1-Epha3 with the
# updated sending/target cluster combinations (ECs/FAPs/sCs x ECs/FAPs/sCs),
# expanding the original 3 data rows into the full 3x3 = 9 row cross table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = the various expression / specificity metrics.
$data = @(
    @("ECs",  "Efna1", "Epha3", "ECs",  3, 1, 11.905615,          35.716845,          0.8197078149061106, 0.8197078149061106, 1, 0.3333333333333333, 0.06698166666666668, 0.200945,  0.003012576978541733, 0.003012576978541732, 0.7974579353916668, 7.177121418525,     0.002469432892316896, 0.002469432892316896),
    @("ECs",  "Efna1", "Epha3", "FAPs", 3, 1, 11.905615,          35.716845,          0.8197078149061106, 0.8197078149061106, 3, 1,                   22.03620333333333,  66.10861,  0.9911034191912899,   0.9911034191912899,   262.3545529483833,  2361.19097653545,  0.8124152180912672,   0.8124152180912672),
    @("ECs",  "Efna1", "Epha3", "sCs",  3, 1, 11.905615,          35.716845,          0.8197078149061106, 0.8197078149061106, 2, 0.6666666666666666, 0.130825,            0.392475,  0.005884003830168287, 0.005884003830168287, 1.557552082375,     14.017968741375,   0.004823163922526432, 0.004823163922526432),
    @("FAPs", "Efna1", "Epha3", "ECs",  3, 1, 1.63134,            4.894019999999999,  0.1123186115768849, 0.1123186115768849, 1, 0.3333333333333333, 0.06698166666666668, 0.200945,  0.003012576978541733, 0.003012576978541732, 0.1092698721,       0.9834288488999999, 0.0003383684634982943, 0.0003383684634982943),
    @("FAPs", "Efna1", "Epha3", "FAPs", 3, 1, 1.63134,            4.894019999999999,  0.1123186115768849, 0.1123186115768849, 3, 1,                   22.03620333333333,  66.10861,  0.9911034191912899,   0.9911034191912899,   35.94853994579999,  323.5368595122,    0.111319359972669,    0.111319359972669),
    @("FAPs", "Efna1", "Epha3", "sCs",  3, 1, 1.63134,            4.894019999999999,  0.1123186115768849, 0.1123186115768849, 2, 0.6666666666666666, 0.130825,            0.392475,  0.005884003830168287, 0.005884003830168287, 0.2134200555,       1.9207804995,      0.0006608831407175748, 0.0006608831407175748),
    @("sCs",  "Efna1", "Epha3", "ECs",  3, 1, 0.9872630000000001, 2.961789,           0.0679735735170045, 0.0679735735170045, 1, 0.3333333333333333, 0.06698166666666668, 0.200945,  0.003012576978541733, 0.003012576978541732, 0.06612852117833334, 0.5951566906050001, 0.0002047756227265418, 0.0002047756227265417),
    @("sCs",  "Efna1", "Epha3", "FAPs", 3, 1, 0.9872630000000001, 2.961789,           0.0679735735170045, 0.0679735735170045, 3, 1,                   22.03620333333333,  66.10861,  0.9911034191912899,   0.9911034191912899,   21.75552821147667,  195.79975390329,   0.06736884112735367,  0.06736884112735367),
    @("sCs",  "Efna1", "Epha3", "sCs",  3, 1, 0.9872630000000001, 2.961789,           0.0679735735170045, 0.0679735735170045, 2, 0.6666666666666666, 0.130825,            0.392475,  0.005884003830168287, 0.005884003830168287, 0.129158681975,     1.162428137775,    0.0003999567669242801, 0.0003999567669242801)
)

for ($rIdx = 0; $rIdx -lt $data.Length; $rIdx++) {
    $rowVals = $data[$rIdx]
    $excelRow = $rIdx + 2
    for ($cIdx = 0; $cIdx -lt $rowVals.Length; $cIdx++) {
        $ws.Cells.Item($excelRow, $cIdx + 1).Value = $rowVals[$cIdx]
    }
}
